$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049992953503481
$ws.Range("D2").Value = 1.054902367190466
$ws.Range("E2").Value = 1.046887125505166
$ws.Range("F2").Value = 1.063816585448425
$ws.Range("I2").Value = 1.041689313359631
$ws.Range("J2").Value = 1.055028462698303
$ws.Range("K2").Value = 1.057644299060754
$ws.Range("L2").Value = 1.049651304281252
$ws.Range("M2").Value = 1.066534201493332
$ws.Range("N2").Value = 1.05652672412879
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.051810594787204
$ws.Range("D3").Value = 1.056365186997686
$ws.Range("E3").Value = 1.048482215375519
$ws.Range("F3").Value = 1.065519848265398
$ws.Range("I3").Value = 1.042202233673633
$ws.Range("J3").Value = 1.056491365792578
$ws.Range("K3").Value = 1.058918796350207
$ws.Range("L3").Value = 1.051056094976843
$ws.Range("M3").Value = 1.068050322900012
$ws.Range("N3").Value = 1.057991704713257
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.052982976709455
$ws.Range("D4").Value = 1.057308392585577
$ws.Range("E4").Value = 1.049510556155923
$ws.Range("F4").Value = 1.066618930886289
$ws.Range("I4").Value = 1.042530968330674
$ws.Range("J4").Value = 1.057433958837811
$ws.Range("K4").Value = 1.059739638572331
$ws.Range("L4").Value = 1.051960827426759
$ws.Range("M4").Value = 1.06902785220985
$ws.Range("N4").Value = 1.058935636348731
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053474966510024
$ws.Range("D5").Value = 1.057704131509535
$ws.Range("E5").Value = 1.049941980200421
$ws.Range("F5").Value = 1.067080272880223
$ws.Range("I5").Value = 1.042668418255341
$ws.Range("J5").Value = 1.057829282165608
$ws.Range("K5").Value = 1.06008381392671
$ws.Range("L5").Value = 1.052340172422737
$ws.Range("M5").Value = 1.069437981902298
$ws.Range("N5").Value = 1.059331521080998
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053557522734263
$ws.Range("D6").Value = 1.057770532207344
$ws.Range("E6").Value = 1.050014366462731
$ws.Range("F6").Value = 1.067157692920688
$ws.Range("I6").Value = 1.042691452892554
$ws.Range("J6").Value = 1.057895603855211
$ws.Range("K6").Value = 1.06014154969502
$ws.Range("L6").Value = 1.052403807660624
$ws.Range("M6").Value = 1.069506796614362
$ws.Range("N6").Value = 1.059397936955007
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.052989554131429
$ws.Range("D7").Value = 1.057313683532571
$ws.Range("E7").Value = 1.049516324341157
$ws.Range("F7").Value = 1.066625098137787
$ws.Range("I7").Value = 1.042532807882335
$ws.Range("J7").Value = 1.057439244850797
$ws.Range("K7").Value = 1.059744241004969
$ws.Range("L7").Value = 1.051965900179055
$ws.Range("M7").Value = 1.069033335598849
$ws.Range("N7").Value = 1.058940929868462
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.050608022597397
$ws.Range("D8").Value = 1.055397433395021
$ws.Range("E8").Value = 1.047426988227253
$ws.Range("F8").Value = 1.06439285070101
$ws.Range("I8").Value = 1.041863315023345
$ws.Range("J8").Value = 1.055523695923184
$ws.Range("K8").Value = 1.058075825068036
$ws.Range("L8").Value = 1.050126951580683
$ws.Range("M8").Value = 1.067047316079682
$ws.Range("N8").Value = 1.057022660641654
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.046381852631049
$ws.Range("D9").Value = 1.051994564360256
$ws.Range("E9").Value = 1.043715538745481
$ws.Range("F9").Value = 1.060435322601125
$ws.Range("I9").Value = 1.040659109259485
$ws.Range("J9").Value = 1.052116892685775
$ws.Range("K9").Value = 1.055105833660423
$ws.Range("L9").Value = 1.046853151173782
$ws.Range("M9").Value = 1.063520196789934
$ws.Range("N9").Value = 1.053611019352884
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043543291336352
$ws.Range("D10").Value = 1.049707478814887
$ws.Range("E10").Value = 1.041220163303661
$ws.Range("F10").Value = 1.057779802328461
$ws.Range("I10").Value = 1.039839456996485
$ws.Range("J10").Value = 1.049823616680121
$ws.Range("K10").Value = 1.053104806604832
$ws.Range("L10").Value = 1.044647213946306
$ws.Range("M10").Value = 1.061149347409864
$ws.Range("N10").Value = 1.051314486632259
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.042308884063246
$ws.Range("D11").Value = 1.048712554836187
$ws.Range("E11").Value = 1.040134402077852
$ws.Range("F11").Value = 1.05662563262117
$ws.Range("I11").Value = 1.039480453683239
$ws.Range("J11").Value = 1.048825145592346
$ws.Range("K11").Value = 1.052233161888467
$ws.Range("L11").Value = 1.043686246082808
$ws.Range("M11").Value = 1.060117920712647
$ws.Range("N11").Value = 1.050314597601013
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041849551787503
$ws.Range("D12").Value = 1.048342287065844
$ws.Range("E12").Value = 1.039730293117656
$ws.Range("F12").Value = 1.056196255002694
$ws.Range("I12").Value = 1.039346481885494
$ws.Range("J12").Value = 1.048453428318462
$ws.Range("K12").Value = 1.051908598530738
$ws.Range("L12").Value = 1.043328411849778
$ws.Range("M12").Value = 1.059734058306078
$ws.Range("N12").Value = 1.049942352445963
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04194811753659
$ws.Range("D13").Value = 1.048421743113088
$ws.Range("E13").Value = 1.039817012798742
$ws.Range("F13").Value = 1.056288388490169
$ws.Range("I13").Value = 1.03937524756227
$ws.Range("J13").Value = 1.048533201288447
$ws.Range("K13").Value = 1.051978254746195
$ws.Range("L13").Value = 1.043405209001855
$ws.Range("M13").Value = 1.059816432089433
$ws.Range("N13").Value = 1.050022238702715
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.042270932347084
$ws.Range("D14").Value = 1.048681962946151
$ws.Range("E14").Value = 1.040101014926342
$ws.Range("F14").Value = 1.056590153867526
$ws.Range("I14").Value = 1.03946939226605
$ws.Range("J14").Value = 1.048794436587825
$ws.Range("K14").Value = 1.052206349708042
$ws.Range("L14").Value = 1.043656685620169
$ws.Range("M14").Value = 1.060086205823836
$ws.Range("N14").Value = 1.050283844986183
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.042469720045859
$ws.Range("D15").Value = 1.048842198536193
$ws.Range("E15").Value = 1.040275890203774
$ws.Range("F15").Value = 1.056775992416121
$ws.Range("I15").Value = 1.03952731527316
$ws.Range("J15").Value = 1.048955280302872
$ws.Range("K15").Value = 1.052346780608356
$ws.Range("L15").Value = 1.043811510418774
$ws.Range("M15").Value = 1.060252323100974
$ws.Range("N15").Value = 1.050444917117755
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043625101400044
$ws.Range("D16").Value = 1.04977341022149
$ws.Range("E16").Value = 1.041292109398716
$ws.Range("F16").Value = 1.057856308241187
$ws.Range("I16").Value = 1.03986319606683
$ws.Range("J16").Value = 1.049889765089245
$ws.Range("K16").Value = 1.053162544115784
$ws.Range("L16").Value = 1.044710866752667
$ws.Range("M16").Value = 1.061217696401176
$ws.Range("N16").Value = 1.051380728979711
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044348408445681
$ws.Range("D17").Value = 1.050356290666393
$ws.Range("E17").Value = 1.041928137594168
$ws.Range("F17").Value = 1.058532793512907
$ws.Range("I17").Value = 1.040072785187775
$ws.Range("J17").Value = 1.050474465419278
$ws.Range("K17").Value = 1.053672850821651
$ws.Range("L17").Value = 1.045273447588463
$ws.Range("M17").Value = 1.061821942607935
$ws.Range("N17").Value = 1.051966259651279
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044769792255122
$ws.Range("D18").Value = 1.050695832246376
$ws.Range("E18").Value = 1.042298617204287
$ws.Range("F18").Value = 1.058926961240092
$ws.Range("I18").Value = 1.040194641025808
$ws.Range("J18").Value = 1.050814985012328
$ws.Range("K18").Value = 1.053970004760855
$ws.Range("L18").Value = 1.04560103485694
$ws.Range("M18").Value = 1.062173924116486
$ws.Range("N18").Value = 1.052307262821209
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044913387622789
$ws.Range("D19").Value = 1.050811532582589
$ws.Range("E19").Value = 1.042424856159951
$ws.Range("F19").Value = 1.059061292505458
$ws.Range("I19").Value = 1.040236124129081
$ws.Range("J19").Value = 1.050931004716321
$ws.Range("K19").Value = 1.054071242374287
$ws.Range("L19").Value = 1.045712639816414
$ws.Range("M19").Value = 1.062293862383979
$ws.Range("N19").Value = 1.052423447286489
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044270857257484
$ws.Range("D20").Value = 1.050293799029789
$ws.Range("E20").Value = 1.041859950157993
$ws.Range("F20").Value = 1.058460256081036
$ws.Range("I20").Value = 1.0400503390625
$ws.Range("J20").Value = 1.050411787175139
$ws.Range("K20").Value = 1.053618151508013
$ws.Range("L20").Value = 1.045213145700878
$ws.Range("M20").Value = 1.061757160947045
$ws.Range("N20").Value = 1.051903492396843
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042175894214613
$ws.Range("D21").Value = 1.048605354396436
$ws.Range("E21").Value = 1.040017405874354
$ws.Range("F21").Value = 1.056501310074427
$ws.Range("I21").Value = 1.039441686229717
$ws.Range("J21").Value = 1.048717532665202
$ws.Range("K21").Value = 1.052139203525276
$ws.Range("L21").Value = 1.043582656692264
$ws.Range("M21").Value = 1.060006784887295
$ws.Range("N21").Value = 1.050206831851169
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040853960638404
$ws.Range("D22").Value = 1.047539653056212
$ws.Range("E22").Value = 1.038854234739366
$ws.Range("F22").Value = 1.055265770449669
$ws.Range("I22").Value = 1.039055399524649
$ws.Range("J22").Value = 1.04764741412573
$ws.Range("K22").Value = 1.051204718245401
$ws.Range("L22").Value = 1.042552356244694
$ws.Range("M22").Value = 1.05890193659296
$ws.Range("N22").Value = 1.049135193620626
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.041555200529603
$ws.Range("D23").Value = 1.048104996901821
$ws.Range("E23").Value = 1.039471305298982
$ws.Range("F23").Value = 1.055921127126515
$ws.Range("I23").Value = 1.039260521503619
$ws.Range("J23").Value = 1.048215172648875
$ws.Range("K23").Value = 1.051700549284131
$ws.Range("L23").Value = 1.043099032394724
$ws.Range("M23").Value = 1.059488053046128
$ws.Range("N23").Value = 1.049703758425997
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044305900873333
$ws.Range("D24").Value = 1.050322037678088
$ws.Range("E24").Value = 1.041890762684857
$ws.Range("F24").Value = 1.058493033905123
$ws.Range("I24").Value = 1.040060482717456
$ws.Range("J24").Value = 1.050440110401241
$ws.Range("K24").Value = 1.053642869315885
$ws.Range("L24").Value = 1.045240395247649
$ws.Range("M24").Value = 1.061786434427797
$ws.Range("N24").Value = 1.051931855845175
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047478050959774
$ws.Range("D25").Value = 1.052877482976442
$ws.Range("E25").Value = 1.044678673753801
$ws.Range("F25").Value = 1.061461389176699
$ws.Range("I25").Value = 1.040973365770205
$ws.Range("J25").Value = 1.053001450434197
$ws.Range("K25").Value = 1.055877292345751
$ws.Range("L25").Value = 1.047703558850671
$ws.Range("M25").Value = 1.064435397273559
$ws.Range("N25").Value = 1.054496833274768
